# Release mCSD 3.9.0 with CP integrated
# Updates the "Metadata" sheet of the mCSD EndpointTypes CodeSystem workbook:
#   - Version bump 3.8.0 -> 3.9.0
#   - Experimental flag now explicitly "false"
#   - Date refreshed to the new publication timestamp
#   - Contact block expanded from a single placeholder row into the
#     three rows FHIR IG Publisher emits for a contact with a url,
#     an email and a resolved name/email label
#   - Jurisdiction changed from "World" to "Global (Whole world)"

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

$ws.Range("B3").Value = "3.9.0"

# Plain Value = "false" gets auto-coerced to a real Boolean by the
# Excel-literal-parsing rules (same as typing FALSE into a cell), but the
# source workbook stores it as literal text "false". Route it through a
# text formula + paste-values-only so it lands as a shared string instead.
$ws.Range("B7").Formula = "=""false"""
$ws.Range("B7").Copy()
$ws.Range("B7").PasteSpecial(-4163)

$ws.Range("B8").Value = "2024-12-02T17:05:26-06:00"
$ws.Range("B10").Value = "null (https://www.ihe.net/ihe_domains/it_infrastructure/)"
$ws.Range("B11").Value = "null (iti@ihe.net)"
$ws.Range("B12").Value = "IHE IT Infrastructure Technical Committee (iti@ihe.net)"
$ws.Range("B13").Value = "Global (Whole world)"
